$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.002365
$ws.Range("H2").Value = 3.007095
$ws.Range("I2").Value = 0.00427078073065012
$ws.Range("J2").Value = 0.00427078073065012
$ws.Range("M2").Value = 159.4836373333333
$ws.Range("N2").Value = 478.450912
$ws.Range("O2").Value = 0.2983285084902258
$ws.Range("P2").Value = 0.2983285084902258
$ws.Range("Q2").Value = 159.8608161356267
$ws.Range("R2").Value = 1438.74734522064
$ws.Range("S2").Value = 0.001274095645463647
$ws.Range("T2").Value = 0.001274095645463647

# Row 3
$ws.Range("G3").Value = 1.002365
$ws.Range("H3").Value = 3.007095
$ws.Range("I3").Value = 0.00427078073065012
$ws.Range("J3").Value = 0.00427078073065012
$ws.Range("O3").Value = 0.3227862111630279
$ws.Range("P3").Value = 0.3227862111630279
$ws.Range("Q3").Value = 172.966597845405
$ws.Range("R3").Value = 1556.699380608645
$ws.Range("S3").Value = 0.00137854913075462
$ws.Range("T3").Value = 0.00137854913075462

# Row 4
$ws.Range("G4").Value = 1.002365
$ws.Range("H4").Value = 3.007095
$ws.Range("I4").Value = 0.00427078073065012
$ws.Range("J4").Value = 0.00427078073065012
$ws.Range("M4").Value = 74.38770566666666
$ws.Range("N4").Value = 223.163117
$ws.Range("O4").Value = 0.1391489036280481
$ws.Range("P4").Value = 0.1391489036280482
$ws.Range("Q4").Value = 74.56363259056832
$ws.Range("R4").Value = 671.0726933151149
$ws.Range("S4").Value = 0.0005942744563057586
$ws.Range("T4").Value = 0.0005942744563057587

# Row 5
$ws.Range("G5").Value = 1.002365
$ws.Range("H5").Value = 3.007095
$ws.Range("I5").Value = 0.00427078073065012
$ws.Range("J5").Value = 0.00427078073065012
$ws.Range("M5").Value = 58.41461433333333
$ws.Range("N5").Value = 175.243843
$ws.Range("O5").Value = 0.1092697975759847
$ws.Range("P5").Value = 0.1092697975759848
$ws.Range("Q5").Value = 58.55276489623166
$ws.Range("R5").Value = 526.9748840660849
$ws.Range("S5").Value = 0.0004666673459295549
$ws.Range("T5").Value = 0.0004666673459295549

# Row 6
$ws.Range("G6").Value = 1.002365
$ws.Range("H6").Value = 3.007095
$ws.Range("I6").Value = 0.00427078073065012
$ws.Range("J6").Value = 0.00427078073065012
$ws.Range("M6").Value = 69.746216
$ws.Range("N6").Value = 209.238648
$ws.Range("O6").Value = 0.1304665791427133
$ws.Range("P6").Value = 0.1304665791427133
$ws.Range("Q6").Value = 69.91116580084
$ws.Range("R6").Value = 629.2004922075599
$ws.Range("S6").Value = 0.0005571941521965389
$ws.Range("T6").Value = 0.000557194152196539

# Row 7
$ws.Range("I7").Value = 0.00206557659722455
$ws.Range("J7").Value = 0.002065576597224551
$ws.Range("M7").Value = 159.4836373333333
$ws.Range("N7").Value = 478.450912
$ws.Range("O7").Value = 0.2983285084902258
$ws.Range("P7").Value = 0.2983285084902258
$ws.Range("Q7").Value = 77.317188928288
$ws.Range("R7").Value = 695.8547003545921
$ws.Range("S7").Value = 0.0006162203854223159
$ws.Range("T7").Value = 0.000616220385422316

# Row 8
$ws.Range("I8").Value = 0.00206557659722455
$ws.Range("J8").Value = 0.002065576597224551
$ws.Range("O8").Value = 0.3227862111630279
$ws.Range("P8").Value = 0.3227862111630279
$ws.Range("S8").Value = 0.0006667396436851323
$ws.Range("T8").Value = 0.0006667396436851324

# Row 9
$ws.Range("I9").Value = 0.00206557659722455
$ws.Range("J9").Value = 0.002065576597224551
$ws.Range("M9").Value = 74.38770566666666
$ws.Range("N9").Value = 223.163117
$ws.Range("O9").Value = 0.1391489036280481
$ws.Range("P9").Value = 0.1391489036280482
$ws.Range("Q9").Value = 36.06293654408299
$ws.Range("R9").Value = 324.566428896747
$ws.Range("S9").Value = 0.0002874227188635505
$ws.Range("T9").Value = 0.0002874227188635506

# Row 10
$ws.Range("I10").Value = 0.00206557659722455
$ws.Range("J10").Value = 0.002065576597224551
$ws.Range("M10").Value = 58.41461433333333
$ws.Range("N10").Value = 175.243843
$ws.Range("O10").Value = 0.1092697975759847
$ws.Range("P10").Value = 0.1092697975759848
$ws.Range("Q10").Value = 28.319229784957
$ws.Range("R10").Value = 254.873068064613
$ws.Range("S10").Value = 0.000225705136656418
$ws.Range("T10").Value = 0.000225705136656418

# Row 11
$ws.Range("I11").Value = 0.00206557659722455
$ws.Range("J11").Value = 0.002065576597224551
$ws.Range("M11").Value = 69.746216
$ws.Range("N11").Value = 209.238648
$ws.Range("O11").Value = 0.1304665791427133
$ws.Range("P11").Value = 0.1304665791427133
$ws.Range("Q11").Value = 33.812756278152
$ws.Range("R11").Value = 304.314806503368
$ws.Range("S11").Value = 0.0002694887125971332
$ws.Range("T11").Value = 0.0002694887125971333

# Row 12
$ws.Range("G12").Value = 135.9134216666667
$ws.Range("H12").Value = 407.740265
$ws.Range("I12").Value = 0.5790868818152315
$ws.Range("J12").Value = 0.5790868818152316
$ws.Range("M12").Value = 159.4836373333333
$ws.Range("N12").Value = 478.450912
$ws.Range("O12").Value = 0.2983285084902258
$ws.Range("P12").Value = 0.2983285084902258
$ws.Range("Q12").Value = 21675.96684981908
$ws.Range("R12").Value = 195083.7016483717
$ws.Range("S12").Value = 0.1727581257381937
$ws.Range("T12").Value = 0.1727581257381937

# Row 13
$ws.Range("G13").Value = 135.9134216666667
$ws.Range("H13").Value = 407.740265
$ws.Range("I13").Value = 0.5790868818152315
$ws.Range("J13").Value = 0.5790868818152316
$ws.Range("O13").Value = 0.3227862111630279
$ws.Range("P13").Value = 0.3227862111630279
$ws.Range("Q13").Value = 23453.01576492723
$ws.Range("R13").Value = 211077.1418843451
$ws.Range("S13").Value = 0.1869212605153507
$ws.Range("T13").Value = 0.1869212605153507

# Row 14
$ws.Range("G14").Value = 135.9134216666667
$ws.Range("H14").Value = 407.740265
$ws.Range("I14").Value = 0.5790868818152315
$ws.Range("J14").Value = 0.5790868818152316
$ws.Range("M14").Value = 74.38770566666666
$ws.Range("N14").Value = 223.163117
$ws.Range("O14").Value = 0.1391489036280481
$ws.Range("P14").Value = 0.1391489036280482
$ws.Range("Q14").Value = 10110.28760708955
$ws.Range("R14").Value = 90992.58846380601
$ws.Range("S14").Value = 0.08057930470997454
$ws.Range("T14").Value = 0.08057930470997457

# Row 15
$ws.Range("G15").Value = 135.9134216666667
$ws.Range("H15").Value = 407.740265
$ws.Range("I15").Value = 0.5790868818152315
$ws.Range("J15").Value = 0.5790868818152316
$ws.Range("M15").Value = 58.41461433333333
$ws.Range("N15").Value = 175.243843
$ws.Range("O15").Value = 0.1092697975759847
$ws.Range("P15").Value = 0.1092697975759848
$ws.Range("Q15").Value = 7939.330109382044
$ws.Range("R15").Value = 71453.9709844384
$ws.Range("S15").Value = 0.06327670635485855
$ws.Range("T15").Value = 0.06327670635485858

# Row 16
$ws.Range("G16").Value = 135.9134216666667
$ws.Range("H16").Value = 407.740265
$ws.Range("I16").Value = 0.5790868818152315
$ws.Range("J16").Value = 0.5790868818152316
$ws.Range("M16").Value = 69.746216
$ws.Range("N16").Value = 209.238648
$ws.Range("O16").Value = 0.1304665791427133
$ws.Range("P16").Value = 0.1304665791427133
$ws.Range("Q16").Value = 9479.446864862413
$ws.Range("R16").Value = 85315.02178376172
$ws.Range("S16").Value = 0.07555148449685398
$ws.Range("T16").Value = 0.075551484496854

# Row 17
$ws.Range("G17").Value = 0.06627866666666667
$ws.Range("H17").Value = 0.198836
$ws.Range("I17").Value = 0.0002823937911371431
$ws.Range("J17").Value = 0.0002823937911371431
$ws.Range("M17").Value = 159.4836373333333
$ws.Range("N17").Value = 478.450912
$ws.Range("O17").Value = 0.2983285084902258
$ws.Range("P17").Value = 0.2983285084902258
$ws.Range("Q17").Value = 10.57036283760356
$ws.Range("R17").Value = 95.13326553843201
$ws.Range("S17").Value = 0.00008424611851684424
$ws.Range("T17").Value = 0.00008424611851684425

# Row 18
$ws.Range("G18").Value = 0.06627866666666667
$ws.Range("H18").Value = 0.198836
$ws.Range("I18").Value = 0.0002823937911371431
$ws.Range("J18").Value = 0.0002823937911371431
$ws.Range("O18").Value = 0.3227862111630279
$ws.Range("P18").Value = 0.3227862111630279
$ws.Range("Q18").Value = 11.436947103164
$ws.Range("R18").Value = 102.932523928476
$ws.Range("S18").Value = 0.00009115282189712186
$ws.Range("T18").Value = 0.00009115282189712188

# Row 19
$ws.Range("G19").Value = 0.06627866666666667
$ws.Range("H19").Value = 0.198836
$ws.Range("I19").Value = 0.0002823937911371431
$ws.Range("J19").Value = 0.0002823937911371431
$ws.Range("M19").Value = 74.38770566666666
$ws.Range("N19").Value = 223.163117
$ws.Range("O19").Value = 0.1391489036280481
$ws.Range("P19").Value = 0.1391489036280482
$ws.Range("Q19").Value = 4.930317947979111
$ws.Range("R19").Value = 44.37286153181201
$ws.Range("S19").Value = 0.00003929478642810147
$ws.Range("T19").Value = 0.00003929478642810149

# Row 20
$ws.Range("G20").Value = 0.06627866666666667
$ws.Range("H20").Value = 0.198836
$ws.Range("I20").Value = 0.0002823937911371431
$ws.Range("J20").Value = 0.0002823937911371431
$ws.Range("M20").Value = 58.41461433333333
$ws.Range("N20").Value = 175.243843
$ws.Range("O20").Value = 0.1092697975759847
$ws.Range("P20").Value = 0.1092697975759848
$ws.Range("Q20").Value = 3.871642751860889
$ws.Range("R20").Value = 34.844784766748
$ws.Range("S20").Value = 0.00003085711239427054
$ws.Range("T20").Value = 0.00003085711239427055

# Row 21
$ws.Range("G21").Value = 0.06627866666666667
$ws.Range("H21").Value = 0.198836
$ws.Range("I21").Value = 0.0002823937911371431
$ws.Range("J21").Value = 0.0002823937911371431
$ws.Range("M21").Value = 69.746216
$ws.Range("N21").Value = 209.238648
$ws.Range("O21").Value = 0.1304665791427133
$ws.Range("P21").Value = 0.1304665791427133
$ws.Range("Q21").Value = 4.622686201525333
$ws.Range("R21").Value = 41.60417581372801
$ws.Range("S21").Value = 0.00003684295190080493
$ws.Range("T21").Value = 0.00003684295190080495

# Row 22
$ws.Range("G22").Value = 97.236126
$ws.Range("H22").Value = 291.708378
$ws.Range("I22").Value = 0.4142943670657566
$ws.Range("J22").Value = 0.4142943670657567
$ws.Range("M22").Value = 159.4836373333333
$ws.Range("N22").Value = 478.450912
$ws.Range("O22").Value = 0.2983285084902258
$ws.Range("P22").Value = 0.2983285084902258
$ws.Range("Q22").Value = 15507.5710546823
$ws.Range("R22").Value = 139568.1394921407
$ws.Range("S22").Value = 0.1235958206026293
$ws.Range("T22").Value = 0.1235958206026293

# Row 23
$ws.Range("G23").Value = 97.236126
$ws.Range("H23").Value = 291.708378
$ws.Range("I23").Value = 0.4142943670657566
$ws.Range("J23").Value = 0.4142943670657567
$ws.Range("O23").Value = 0.3227862111630279
$ws.Range("P23").Value = 0.3227862111630279
$ws.Range("Q23").Value = 16778.91975666262
$ws.Range("R23").Value = 151010.2778099636
$ws.Range("S23").Value = 0.1337285090513403
$ws.Range("T23").Value = 0.1337285090513403

# Row 24
$ws.Range("G24").Value = 97.236126
$ws.Range("H24").Value = 291.708378
$ws.Range("I24").Value = 0.4142943670657566
$ws.Range("J24").Value = 0.4142943670657567
$ws.Range("M24").Value = 74.38770566666666
$ws.Range("N24").Value = 223.163117
$ws.Range("O24").Value = 0.1391489036280481
$ws.Range("P24").Value = 0.1391489036280482
$ws.Range("Q24").Value = 7233.172321054913
$ws.Range("R24").Value = 65098.55088949422
$ws.Range("S24").Value = 0.05764860695647617
$ws.Range("T24").Value = 0.05764860695647619

# Row 25
$ws.Range("G25").Value = 97.236126
$ws.Range("H25").Value = 291.708378
$ws.Range("I25").Value = 0.4142943670657566
$ws.Range("J25").Value = 0.4142943670657567
$ws.Range("M25").Value = 58.41461433333333
$ws.Range("N25").Value = 175.243843
$ws.Range("O25").Value = 0.1092697975759847
$ws.Range("P25").Value = 0.1092697975759848
$ws.Range("Q25").Value = 5680.010799557406
$ws.Range("R25").Value = 51120.09719601665
$ws.Range("S25").Value = 0.04526986162614595
$ws.Range("T25").Value = 0.04526986162614596

# Row 26
$ws.Range("G26").Value = 97.236126
$ws.Range("H26").Value = 291.708378
$ws.Range("I26").Value = 0.4142943670657566
$ws.Range("J26").Value = 0.4142943670657567
$ws.Range("M26").Value = 69.746216
$ws.Range("N26").Value = 209.238648
$ws.Range("O26").Value = 0.1304665791427133
$ws.Range("P26").Value = 0.1304665791427133
$ws.Range("Q26").Value = 6781.851846999216
$ws.Range("R26").Value = 61036.66662299295
$ws.Range("S26").Value = 0.05405156882916486
$ws.Range("T26").Value = 0.05405156882916488
